$wb = $excel.ActiveWorkbook

# --- Sheet3: add the new "percent_of_control" column (D) ---
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("D1").Value = "percent_of_control"

$ws3.Range("D2").Formula = "=(B2/`$I`$1)*100"
$ws3.Range("D3:D6").Formula = "=(B3/`$I`$1)*100"

# --- Selection / active-sheet bookkeeping ---
# Sheet1 had tabSelected; move the selection there first (matches its own
# selection rectangle) so activating Sheet3 afterwards leaves Sheet1 without
# tabSelected, like the target file.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("B9:H9").Select() | Out-Null

# Sheet3 becomes the active tab, with D4 selected.
$ws3.Activate() | Out-Null
$ws3.Range("D4").Select() | Out-Null
